$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Modelo" column
$ws.Range("F1").Value = "Modelo"

# Match formatting of the other header cells (bold, border, centered)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the numeric prediction values
$ws.Range("B2").Value = 0.05027256464626385
$ws.Range("C2").Value = 0.9985215696412008
$ws.Range("D2").Value = 0.1657370521520764

# Add the model name value
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
